$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Average overlap) values
$ws.Range("C2").Value = 85.9731668718597
$ws.Range("C3").Value = 33.06992194236347
$ws.Range("C4").Value = 31.49790873184837

# Remove column D (Lowest overlap) entirely
$ws.Range("D1:D4").Delete()
